$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2 = @{ "A"="ECs"; "B"="Icam2"; "C"="Itgb2"; "D"="ECs"; "E"="3"; "F"="1"; "G"="23.76874"; "H"="71.30622000000001"; "I"="0.9555483638834226"; "J"="0.9555483638834227"; "K"="2"; "L"="0.6666666666666666"; "M"="83.95844533333333"; "N"="251.875336"; "O"="0.9979754487867319"; "P"="0.9979754487867319"; "Q"="1995.586457932214"; "R"="17960.27812138992"; "S"="0.953613807283986"; "T"="0.9536138072839861" }
  3 = @{ "A"="ECs"; "B"="Icam2"; "C"="Itgb2"; "D"="FAPs"; "E"="3"; "F"="1"; "G"="23.76874"; "H"="71.30622000000001"; "I"="0.9555483638834226"; "J"="0.9555483638834227"; "K"="3"; "L"="1"; "M"="0.170323"; "N"="0.510969"; "O"="0.002024551213268089"; "P"="0.00202455121326809"; "Q"="4.048363103020001"; "R"="36.43526792718001"; "S"="0.001934556599436521"; "T"="0.001934556599436522" }
  4 = @{ "A"="FAPs"; "B"="Icam2"; "C"="Itgb2"; "D"="ECs"; "E"="3"; "F"="1"; "G"="0.7500946666666666"; "H"="2.250284"; "I"="0.03015522621270687"; "J"="0.03015522621270688"; "K"="2"; "L"="0.6666666666666666"; "M"="83.95844533333333"; "N"="251.875336"; "O"="0.9979754487867319"; "P"="0.9979754487867319"; "Q"="62.97678206615821"; "R"="566.7910385954239"; "S"="0.03009417541289157"; "T"="0.03009417541289157" }
  5 = @{ "A"="FAPs"; "B"="Icam2"; "C"="Itgb2"; "D"="FAPs"; "E"="3"; "F"="1"; "G"="0.7500946666666666"; "H"="2.250284"; "I"="0.03015522621270687"; "J"="0.03015522621270688"; "K"="3"; "L"="1"; "M"="0.170323"; "N"="0.510969"; "O"="0.002024551213268089"; "P"="0.00202455121326809"; "Q"="0.1277583739106667"; "R"="1.149825365196"; "S"="6.105079981530939E-05"; "T"="6.105079981530941E-05" }
  6 = @{ "A"="sCs"; "B"="Icam2"; "C"="Itgb2"; "D"="ECs"; "E"="3"; "F"="1"; "G"="0.3556153333333333"; "H"="1.066846"; "I"="0.01429640990387057"; "J"="0.01429640990387057"; "K"="2"; "L"="0.6666666666666666"; "M"="83.95844533333333"; "N"="251.875336"; "O"="0.9979754487867319"; "P"="0.9979754487867319"; "Q"="29.85691052336178"; "R"="268.712194710256"; "S"="0.01426746608985431"; "T"="0.01426746608985431" }
  7 = @{ "A"="sCs"; "B"="Icam2"; "C"="Itgb2"; "D"="FAPs"; "E"="3"; "F"="1"; "G"="0.3556153333333333"; "H"="1.066846"; "I"="0.01429640990387057"; "J"="0.01429640990387057"; "K"="3"; "L"="1"; "M"="0.170323"; "N"="0.510969"; "O"="0.002024551213268089"; "P"="0.00202455121326809"; "Q"="0.06056947041933333"; "R"="0.545125233774"; "S"="2.894381401625909E-05"; "T"="2.89438140162591E-05" }
}

foreach ($r in $data.Keys) {
  $row = $data[$r]
  foreach ($c in $row.Keys) {
    $addr = "$c$r"
    $val = $row[$c]
    if ($c -in @("A","B","C","D")) {
      $ws.Range($addr).Value = $val
    } else {
      $ws.Range($addr).Value = [double]$val
    }
  }
}

